# Update "想去人数" (F column) counts for two worksheets: "展览" and "全部类型".
# The two sheets list the same events but the "全部类型" sheet rows are offset by 1
# relative to the "展览" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value, for the "展览" sheet.
$exhibitUpdates = @{
    5  = 1220
    7  = 13960
    8  = 15092
    10 = 23
    21 = 1171
    27 = 5474
    28 = 63
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new F value, for the "全部类型" sheet (rows offset by +1).
$allTypeUpdates = @{
    6  = 1220
    8  = 13960
    9  = 15092
    11 = 23
    22 = 1171
    29 = 5474
    30 = 63
}

foreach ($row in $allTypeUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypeUpdates[$row]
}
